$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D: "canonical SMILES"
$ws.Range("D2").Value = "canonical SMILES"

# Duplicate the "canonical isomeric SMILES" values (column C) into the new
# "canonical SMILES" column D, row by row (Text avoids COM Variant issues).
$ws.Range("D3").Value = $ws.Range("C3").Text
$ws.Range("D4").Value = $ws.Range("C4").Text
$ws.Range("D5").Value = $ws.Range("C5").Text
$ws.Range("D6").Value = $ws.Range("C6").Text
$ws.Range("D7").Value = $ws.Range("C7").Text
$ws.Range("D8").Value = $ws.Range("C8").Text

# Set column D width to match the target layout (~36.86 characters wide)
$ws.Columns.Item(4).ColumnWidth = 36
